$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 4 ("Solid State Relay"): Manufacturer Part Number cell ---
# 10PCV2415 -> 4D2425, paragraph becomes a (de-bolded) Heading1-styled
# paragraph with explicit zero spacing, and the stray _GoBack bookmark is
# removed from here (it moves to the Price cell below).
$cellPart = $t.Cell(4, 2)
$rngPart = $cellPart.Range
$xmlPart = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="007A4D2B" w:rsidRPr="007A4D2B" w:rsidRDefault="007A4D2B" w:rsidP="007A4D2B">
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="007A4D2B">
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>4D2425</w:t>
  </w:r>
</w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$rngPart.InsertXML($xmlPart)

# --- Row 4 ("Solid State Relay"): Price cell ---
# "0.6" + "8" (two runs) -> a single "61.11" run, and the _GoBack bookmark
# now lives inside the hyperlink, right after the run.
$cellPrice = $t.Cell(4, 3)
$rngPrice = $cellPrice.Range
$xmlPrice = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p w:rsidR="00CA1A3E" w:rsidRDefault="007C1687" w:rsidP="00A775D6">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:hyperlink r:id="rId6" w:history="1">
    <w:r w:rsidR="00A775D6" w:rsidRPr="00A775D6">
      <w:rPr>
        <w:bCs/>
        <w:color w:val="000000"/>
      </w:rPr>
      <w:t>61.11</w:t>
    </w:r>
    <w:bookmarkStart w:id="0" w:name="_GoBack"/>
    <w:bookmarkEnd w:id="0"/>
  </w:hyperlink>
</w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$rngPrice.InsertXML($xmlPrice)

Write-Output "Updated parts list: 10PCV2415 -> 4D2425, 0.68 -> 61.11"
